$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320, shifting existing rows 320-429 down to 321-430.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new data record.
$ws.Cells.Item(320, 1).Value = 4
$ws.Cells.Item(320, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(320, 3).Value = "Los Lagos"
$ws.Cells.Item(320, 4).Value = 44809
$ws.Cells.Item(320, 5).Value = 10
$ws.Cells.Item(320, 6).Value = 100114013
$ws.Cells.Item(320, 7).Value = "Zanahoria"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 250
$ws.Cells.Item(320, 11).Value = 10000
$ws.Cells.Item(320, 12).Value = 10000
$ws.Cells.Item(320, 13).Value = 10000
$ws.Cells.Item(320, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(320, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(320, 16).Value = 500
$ws.Cells.Item(320, 17).Value = 20
$ws.Cells.Item(320, 18).Value = "Hortaliza"
